{"js": "const pairs = [\n  [\"91\u00d726=\", \"90\u00d748=\"],\n  [\"14\u00d756=\", \"23\u00d764=\"],\n  [\"80\u00d732=\", \"46\u00d791=\"],\n  [\"93\u00d735=\", \"24\u00d730=\"],\n  [\"67\u00d784=\", \"95\u00d750=\"],\n  [\"65\u00d795=\", \"59\u00d743=\"],\n  [\"56\u00d727=\", \"70\u00d743=\"],\n  [\"24\u00d762=\", \"26\u00d799=\"],\n  [\"42\u00d723=\", \"63\u00d729=\"],\n  [\"65\u00d766=\", \"41\u00d791=\"],\n  [\"89\u00d786=\", \"81\u00d752=\"],\n  [\"75\u00d793=\", \"28\u00d775=\"],\n  [\"18\u00d758=\", \"97\u00d772=\"],\n  [\"54\u00d728=\", \"51\u00d781=\"],\n  [\"95\u00d758=\", \"44\u00d793=\"],\n  [\"15\u00d711=\", \"92\u00d732=\"],\n  [\"40\u00d712=\", \"81\u00d751=\"],\n  [\"41\u00d764=\", \"71\u00d783=\"],\n  [\"57\u00d722=\", \"46\u00d790=\"],\n  [\"81\u00d794=\", \"82\u00d771=\"],\n  [\"47\u00d766=\", \"70\u00d768=\"],\n  [\"89\u00d749=\", \"28\u00d725=\"],\n  [\"80\u00d738=\", \"68\u00d721=\"],\n  [\"68\u00d760=\", \"39\u00d716=\"],\n  [\"94\u00d738=\", \"55\u00d793=\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = [ordered]@{\n    \"91\u00d726=\" = \"90\u00d748=\"\n    \"14\u00d756=\" = \"23\u00d764=\"\n    \"80\u00d732=\" = \"46\u00d791=\"\n    \"93\u00d735=\" = \"24\u00d730=\"\n    \"67\u00d784=\" = \"95\u00d750=\"\n    \"65\u00d795=\" = \"59\u00d743=\"\n    \"56\u00d727=\" = \"70\u00d743=\"\n    \"24\u00d762=\" = \"26\u00d799=\"\n    \"42\u00d723=\" = \"63\u00d729=\"\n    \"65\u00d766=\" = \"41\u00d791=\"\n    \"89\u00d786=\" = \"81\u00d752=\"\n    \"75\u00d793=\" = \"28\u00d775=\"\n    \"18\u00d758=\" = \"97\u00d772=\"\n    \"54\u00d728=\" = \"51\u00d781=\"\n    \"95\u00d758=\" = \"44\u00d793=\"\n    \"15\u00d711=\" = \"92\u00d732=\"\n    \"40\u00d712=\" = \"81\u00d751=\"\n    \"41\u00d764=\" = \"71\u00d783=\"\n    \"57\u00d722=\" = \"46\u00d790=\"\n    \"81\u00d794=\" = \"82\u00d771=\"\n    \"47\u00d766=\" = \"70\u00d768=\"\n    \"89\u00d749=\" = \"28\u00d725=\"\n    \"80\u00d738=\" = \"68\u00d721=\"\n    \"68\u00d760=\" = \"39\u00d716=\"\n    \"94\u00d738=\" = \"55\u00d793=\"\n}\n\nforeach ($key in $pairs.Keys) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $key\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $pairs[$key]\n    $find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2, $true, $false, $false, $false) | Out-Null\n}\n"}
